$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 113: hasMany / cuenta_bancaria / 1 a N / transaccion_bancaria / belongsTo
$ws.Range("A113").Value = "✅"
$ws.Range("B113").Value = "hasMany"
$ws.Range("C113").Value = "cuenta_bancaria"
$ws.Range("D113").Value = "1 a N"
$ws.Range("E113").Value = "transaccion_bancaria"
$ws.Range("F113").Value = "belongsTo"

# Row 114: hasMany / cuenta_bancaria / 1 a N / transaccion_bancaria / belongsTo
$ws.Range("A114").Value = "✅"
$ws.Range("B114").Value = "hasMany"
$ws.Range("C114").Value = "cuenta_bancaria"
$ws.Range("D114").Value = "1 a N"
$ws.Range("E114").Value = "transaccion_bancaria"
$ws.Range("F114").Value = "belongsTo"

# Row 126: belongsTo / concepto_financiero / N a 1 / plan_cuenta / hasMany
$ws.Range("A126").Value = "✅"
$ws.Range("B126").Value = "belongsTo"
$ws.Range("C126").Value = "concepto_financiero"
$ws.Range("D126").Value = "N a 1"
$ws.Range("E126").Value = "plan_cuenta"
$ws.Range("F126").Value = "hasMany"

# Row 140: hasOne / registro_transaccion / 1 a 1 / transaccion_bancaria / belongsTo
$ws.Range("A140").Value = "✅"
$ws.Range("B140").Value = "hasOne"
$ws.Range("C140").Value = "registro_transaccion"
$ws.Range("D140").Value = "1 a 1"
$ws.Range("E140").Value = "transaccion_bancaria"
$ws.Range("F140").Value = "belongsTo"

# Update the saved view state (scroll position + active selection)
$ws.Application.ActiveWindow.ScrollRow = 97
$ws.Range("C119").Select()
